$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(98, 8).Value = 1144.5555
$ws.Cells.Item(98, 9).Value = 1049.7142
$ws.Cells.Item(98, 10).Value = 1476.5
$ws.Cells.Item(98, 11).Value = 1049.7142
$ws.Cells.Item(98, 12).Value = 1476.5
$ws.Cells.Item(98, 13).Value = 448.2858000000001
$ws.Cells.Item(98, 14).Value = -4472.5

$ws.Cells.Item(122, 8).Value = 1144.5555
$ws.Cells.Item(122, 9).Value = 1049.7142
$ws.Cells.Item(122, 10).Value = 1476.5
$ws.Cells.Item(122, 11).Value = 3149.1426
$ws.Cells.Item(122, 12).Value = 4429.5
$ws.Cells.Item(122, 13).Value = -699.1425999999997
$ws.Cells.Item(122, 14).Value = -9329.5

$ws.Cells.Item(132, 8).Value = 5925.7393
$ws.Cells.Item(132, 9).Value = 6383.8423
$ws.Cells.Item(132, 10).Value = 3749.75
$ws.Cells.Item(132, 11).Value = 19151.5269
$ws.Cells.Item(132, 12).Value = 11249.25
$ws.Cells.Item(132, 13).Value = -16621.5269
$ws.Cells.Item(132, 14).Value = -16309.25

$ws.Cells.Item(138, 8).Value = 2640.541
$ws.Cells.Item(138, 9).Value = 4049.1
$ws.Cells.Item(138, 10).Value = 2364.353
$ws.Cells.Item(138, 11).Value = 12147.3
$ws.Cells.Item(138, 12).Value = 7093.059
$ws.Cells.Item(138, 13).Value = -7007.299999999999
$ws.Cells.Item(138, 14).Value = -17373.059

$ws.Cells.Item(141, 8).Value = 4818.6665
$ws.Cells.Item(141, 9).Value = 1809.3529
$ws.Cells.Item(141, 10).Value = 12127
$ws.Cells.Item(141, 11).Value = 5428.0587
$ws.Cells.Item(141, 12).Value = 36381
$ws.Cells.Item(141, 13).Value = -248.0587000000005
$ws.Cells.Item(141, 14).Value = -46741

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 521888.3
$ws.Cells.Item(32, 9).Value = 602833.9
$ws.Cells.Item(32, 11).Value = 602833.9
$ws.Cells.Item(32, 13).Value = -602546.9

$ws.Cells.Item(74, 8).Value = 1646.3125
$ws.Cells.Item(74, 9).Value = 1082.3334
$ws.Cells.Item(74, 11).Value = 1082.3334
$ws.Cells.Item(74, 13).Value = -208.3334

$ws.Cells.Item(77, 8).Value = 1646.3125
$ws.Cells.Item(77, 9).Value = 1082.3334
$ws.Cells.Item(77, 11).Value = 5411.666999999999
$ws.Cells.Item(77, 13).Value = -1043.666999999999

$ws.Cells.Item(107, 8).Value = 29720
$ws.Cells.Item(107, 10).Value = 29720
$ws.Cells.Item(107, 12).Value = 29720
$ws.Cells.Item(107, 14).Value = -37400

$ws.Cells.Item(112, 8).Value = 48043.5
$ws.Cells.Item(112, 10).Value = 48043.5
$ws.Cells.Item(112, 12).Value = 48043.5
$ws.Cells.Item(112, 14).Value = -50997.5

$ws.Cells.Item(141, 8).Value = 100354
$ws.Cells.Item(141, 10).Value = 100354
$ws.Cells.Item(141, 12).Value = 100354
$ws.Cells.Item(141, 14).Value = -110714

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(112, 8).Value = 100469
$ws.Cells.Item(112, 10).Value = 100469
$ws.Cells.Item(112, 12).Value = 100469
$ws.Cells.Item(112, 14).Value = -103423

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 2907.029
$ws.Cells.Item(31, 9).Value = 941.2195
$ws.Cells.Item(31, 10).Value = 5785.5356
$ws.Cells.Item(31, 11).Value = 941.2195
$ws.Cells.Item(31, 12).Value = 5785.5356
$ws.Cells.Item(31, 13).Value = -646.2195
$ws.Cells.Item(31, 14).Value = -6375.5356

$ws.Cells.Item(34, 8).Value = 2907.029
$ws.Cells.Item(34, 9).Value = 941.2195
$ws.Cells.Item(34, 10).Value = 5785.5356
$ws.Cells.Item(34, 11).Value = 941.2195
$ws.Cells.Item(34, 12).Value = 5785.5356
$ws.Cells.Item(34, 13).Value = -739.2195
$ws.Cells.Item(34, 14).Value = -6189.5356

$ws.Cells.Item(119, 8).Value = 0
$ws.Cells.Item(119, 10).Value = 0
$ws.Cells.Item(119, 12).Value = 0
$ws.Cells.Item(119, 14).ClearContents()

$ws.Cells.Item(132, 8).Value = 10418994
$ws.Cells.Item(132, 9).Value = 1122.6666
$ws.Cells.Item(132, 10).Value = 16669717
$ws.Cells.Item(132, 11).Value = 3367.9998
$ws.Cells.Item(132, 12).Value = 50009151
$ws.Cells.Item(132, 13).Value = -837.9998000000001
$ws.Cells.Item(132, 14).Value = -50014211

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(5, 8).Value = 903.7377300000001
$ws.Cells.Item(5, 9).Value = 419.2069
$ws.Cells.Item(5, 11).Value = 1257.6207
$ws.Cells.Item(5, 13).Value = -1145.6207

$ws.Cells.Item(68, 8).Value = 2316.8062
$ws.Cells.Item(68, 10).Value = 1722.0469
$ws.Cells.Item(68, 12).Value = 5166.1407
$ws.Cells.Item(68, 14).Value = -6788.1407

$ws.Cells.Item(71, 8).Value = 2316.8062
$ws.Cells.Item(71, 10).Value = 1722.0469
$ws.Cells.Item(71, 12).Value = 15498.4221
$ws.Cells.Item(71, 14).Value = -23610.4221

$ws.Cells.Item(107, 8).Value = 1712.4412
$ws.Cells.Item(107, 9).Value = 370.57144
$ws.Cells.Item(107, 10).Value = 2060.3333
$ws.Cells.Item(107, 11).Value = 1111.71432
$ws.Cells.Item(107, 12).Value = 6180.999899999999
$ws.Cells.Item(107, 13).Value = 808.28568
$ws.Cells.Item(107, 14).Value = -10020.9999

$ws.Cells.Item(122, 8).Value = 8200.538
$ws.Cells.Item(122, 9).Value = 419
$ws.Cells.Item(122, 10).Value = 50999
$ws.Cells.Item(122, 11).Value = 3771
$ws.Cells.Item(122, 12).Value = 458991
$ws.Cells.Item(122, 13).Value = -1321
$ws.Cells.Item(122, 14).Value = -463891

$ws.Cells.Item(132, 8).Value = 3241.6812
$ws.Cells.Item(132, 9).Value = 2161.1516
$ws.Cells.Item(132, 10).Value = 4232.1665
$ws.Cells.Item(132, 11).Value = 19450.3644
$ws.Cells.Item(132, 12).Value = 38089.4985
$ws.Cells.Item(132, 13).Value = -16920.3644
$ws.Cells.Item(132, 14).Value = -43149.4985

$ws.Cells.Item(135, 8).Value = 903.7377300000001
$ws.Cells.Item(135, 9).Value = 419.2069
$ws.Cells.Item(135, 11).Value = 3772.8621
$ws.Cells.Item(135, 13).Value = -1237.8621

$ws.Cells.Item(137, 8).Value = 7195.879
$ws.Cells.Item(137, 9).Value = 10216.125
$ws.Cells.Item(137, 10).Value = 4353.294
$ws.Cells.Item(137, 11).Value = 30648.375
$ws.Cells.Item(137, 12).Value = 13059.882
$ws.Cells.Item(137, 13).Value = -25548.375
$ws.Cells.Item(137, 14).Value = -23259.882

$ws.Cells.Item(140, 8).Value = 1939.2632
$ws.Cells.Item(140, 9).Value = 1400.4286
$ws.Cells.Item(140, 10).Value = 3448
$ws.Cells.Item(140, 11).Value = 4201.2858
$ws.Cells.Item(140, 12).Value = 10344
$ws.Cells.Item(140, 13).Value = 978.7142000000003
$ws.Cells.Item(140, 14).Value = -20704

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(55, 8).Value = 2945
$ws.Cells.Item(55, 9).Value = 990
$ws.Cells.Item(55, 10).Value = 4900
$ws.Cells.Item(55, 11).Value = 990
$ws.Cells.Item(55, 12).Value = 4900
$ws.Cells.Item(55, 13).Value = -663
$ws.Cells.Item(55, 14).Value = -5554

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 8).Value = 6033.75
$ws.Cells.Item(22, 9).Value = 698.1667
$ws.Cells.Item(22, 10).Value = 14037.125
$ws.Cells.Item(22, 11).Value = 698.1667
$ws.Cells.Item(22, 12).Value = 14037.125
$ws.Cells.Item(22, 13).Value = -403.1667
$ws.Cells.Item(22, 14).Value = -14627.125

$ws.Cells.Item(27, 8).Value = 6033.75
$ws.Cells.Item(27, 9).Value = 698.1667
$ws.Cells.Item(27, 10).Value = 14037.125
$ws.Cells.Item(27, 11).Value = 698.1667
$ws.Cells.Item(27, 12).Value = 14037.125
$ws.Cells.Item(27, 13).Value = -591.1667
$ws.Cells.Item(27, 14).Value = -14251.125

$ws.Cells.Item(45, 8).Value = 0
$ws.Cells.Item(45, 9).Value = 0
$ws.Cells.Item(45, 10).Value = 0
$ws.Cells.Item(45, 11).Value = 0
$ws.Cells.Item(45, 12).Value = 0
$ws.Cells.Item(45, 13).ClearContents()
$ws.Cells.Item(45, 14).ClearContents()

$ws.Cells.Item(93, 8).Value = 11255.272
$ws.Cells.Item(93, 9).Value = 51500.5
$ws.Cells.Item(93, 10).Value = 2311.889
$ws.Cells.Item(93, 11).Value = 51500.5
$ws.Cells.Item(93, 12).Value = 2311.889
$ws.Cells.Item(93, 13).Value = -50252.5
$ws.Cells.Item(93, 14).Value = -4807.889

$ws.Cells.Item(110, 8).Value = 66881.336
$ws.Cells.Item(110, 10).Value = 66881.336
$ws.Cells.Item(110, 12).Value = 66881.336
$ws.Cells.Item(110, 14).Value = -75061.336

$ws.Cells.Item(132, 8).Value = 2639.3408
$ws.Cells.Item(132, 9).Value = 2432.4546
$ws.Cells.Item(132, 10).Value = 3260
$ws.Cells.Item(132, 11).Value = 7297.3638
$ws.Cells.Item(132, 12).Value = 9780
$ws.Cells.Item(132, 13).Value = -4767.3638
$ws.Cells.Item(132, 14).Value = -14840
